# Append a new row (row 69) of sensor-log data to each of the four
# worksheets (ROW35-FE-LIFTER, ROW35-MID-LIFTER, ROW02-FE-LIFTER,
# ROW02-MID-LIFTER). Every sheet currently ends at row 68; this adds one
# more reading to each, matching columns:
#   A time | B total-len(hex) | C ID(hex) | D actual-len(hex)
#   E checksum(hex) | F total-len(dec) | G ID(dec, big int as text)
#   H actual-len(dec) | I checksum(dec)

$wb = $excel.ActiveWorkbook

$rows = @(
    @{ Sheet = 1; A = "2025-03-07 04:42:06"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; D = "0x01,0x90,"; E = "0x d";  F = 400; G = "568631262647113770877196"; H = 400; I = 13  },
    @{ Sheet = 2; A = "2025-03-07 04:29:35"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; D = "0x01,0x90,"; E = "0x e";  F = 400; G = "568631262647113770942732"; H = 400; I = 14  },
    @{ Sheet = 3; A = "2025-03-07 04:51:45"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"; D = "0x01,0x90,"; E = "0xff";  F = 400; G = "568631262647113769959692"; H = 400; I = 255 },
    @{ Sheet = 4; A = "2025-03-07 04:41:15"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x01,0x90,"; E = "0x 3";  F = 400; G = "568631262647113769959692"; H = 400; I = 3   }
)

$newRow = 69

foreach ($entry in $rows) {
    $ws = $wb.Worksheets.Item($entry.Sheet)

    # Text columns (A-E): plain strings, never numeric-looking enough to be
    # auto-coerced by Excel, so a direct .Value assignment keeps them as text.
    $ws.Cells.Item($newRow, 1).Value = $entry.A
    $ws.Cells.Item($newRow, 2).Value = $entry.B
    $ws.Cells.Item($newRow, 3).Value = $entry.C
    $ws.Cells.Item($newRow, 4).Value = $entry.D
    $ws.Cells.Item($newRow, 5).Value = $entry.E

    # Numeric columns (F, H, I): plain integers.
    $ws.Cells.Item($newRow, 6).Value = $entry.F
    $ws.Cells.Item($newRow, 8).Value = $entry.H
    $ws.Cells.Item($newRow, 9).Value = $entry.I

    # Column G holds a ~24-digit id that must stay exact text - format the
    # cell as text first so Excel doesn't round-trip it through a double
    # and mangle it into scientific notation.
    $gCell = $ws.Cells.Item($newRow, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $entry.G
}

Write-Host "Appended row 69 to $($rows.Count) worksheets"
